$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority (E) "low" -> "ht"; Latest Handoff Datetime (H) updated
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$zhcn.Range("H4").Value = "2016-09-05 04:36:05"
$zhcn.Range("H5").Value = "2016-09-05 04:36:05"
$zhcn.Range("H6").Value = "2016-09-05 04:36:05"
$zhcn.Range("H7").Value = "2016-09-05 04:36:05"

# de-de sheet: rows 4-7 -> Priority (E) "low" -> "ht"; Latest Handoff Datetime (H) updated
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

$dede.Range("H4").Value = "2016-09-05 04:36:12"
$dede.Range("H5").Value = "2016-09-05 04:36:12"
$dede.Range("H6").Value = "2016-09-05 04:36:12"
$dede.Range("H7").Value = "2016-09-05 04:36:12"

# Overview sheet: de-de "Latest HO Xliff Generate Date" (G) mirrors the same
# updated handoff timestamp for rows 4-7
$overview.Range("G4").Value = "2016-09-05 04:36:12"
$overview.Range("G5").Value = "2016-09-05 04:36:12"
$overview.Range("G6").Value = "2016-09-05 04:36:12"
$overview.Range("G7").Value = "2016-09-05 04:36:12"
